$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit rotates/swaps species-record data between rows while leaving the
# shared location/meta columns (D, P, S, T, U, V, W, Y, AA, AD, AE, AG, AT,
# AW, AX, AY) untouched:
#   rows 6,7,8  : data cyclically shifts down (old row6 -> row7, old row7 -> row8, old row8 -> row6)
#   rows 17,18  : data is swapped (old row17 <-> old row18)

# --- Row 6 gets old Row 8's data ---
$ws.Range("A6").Value = 131198466
$ws.Range("B6").Value = 79835
$ws.Range("E6").Value = 229821
$ws.Range("F6").Value = "Vedflamlav"
$ws.Range("G6").Value = "Ramboldia elabens"
$ws.Range("H6").Value = "(Fr.) Kantvilas & Elix"
$ws.Range("M6").ClearContents()
$ws.Range("Q6").Value = 466092
$ws.Range("R6").Value = 6789074

# --- Row 7 gets old Row 6's data ---
$ws.Range("A7").Value = 131198844
$ws.Range("B7").Value = 79245
$ws.Range("E7").Value = 6425
$ws.Range("F7").Value = "Garnlav"
$ws.Range("G7").Value = "Alectoria sarmentosa"
$ws.Range("H7").Value = "(Ach.) Ach."
$ws.Range("M7").ClearContents()
$ws.Range("Q7").Value = 466309
$ws.Range("R7").Value = 6789077

# --- Row 8 gets old Row 7's data ---
$ws.Range("A8").Value = 131197802
$ws.Range("B8").Value = 57881
$ws.Range("E8").Value = 100049
$ws.Range("F8").Value = "Spillkråka"
$ws.Range("G8").Value = "Dryocopus martius"
$ws.Range("H8").Value = "(Linnaeus, 1758)"
$ws.Range("M8").Value = "äldre spår"
$ws.Range("Q8").Value = 465938
$ws.Range("R8").Value = 6789021

# --- Row 17 gets old Row 18's data ---
$ws.Range("A17").Value = 131198252
$ws.Range("B17").Value = 79245
$ws.Range("E17").Value = 6425
$ws.Range("F17").Value = "Garnlav"
$ws.Range("G17").Value = "Alectoria sarmentosa"
$ws.Range("H17").Value = "(Ach.) Ach."
$ws.Range("M17").ClearContents()
$ws.Range("Q17").Value = 466111
$ws.Range("R17").Value = 6789063

# --- Row 18 gets old Row 17's data ---
$ws.Range("A18").Value = 131198195
$ws.Range("B18").Value = 57881
$ws.Range("E18").Value = 100049
$ws.Range("F18").Value = "Spillkråka"
$ws.Range("G18").Value = "Dryocopus martius"
$ws.Range("H18").Value = "(Linnaeus, 1758)"
$ws.Range("M18").Value = "färska spår"
$ws.Range("Q18").Value = 466050
$ws.Range("R18").Value = 6788971
